# "using gains for all" - insert two new metric columns (M_TotalTax, M_CorpTax)
# right after M_POP (column E), shifting the existing GFA/IMF/OECD columns
# right by two, and update the row-6 E value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at F:G - everything from old column F onward
# (GFA - Sales .. OECD - Sales + Emp) slides right to H:O. Number formatting
# / style on the old F:G header cells carries over to the new ones.
$ws.Range("F1:G1").EntireColumn.Insert()

# New header row cells for the inserted columns.
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# New data values for rows 2-6, columns F (M_TotalTax) and G (M_CorpTax).
$ws.Range("F2").Value = 14106286460237.92
$ws.Range("G2").Value = 1155021202746.413

$ws.Range("F3").Value = 3207987015.574299
$ws.Range("G3").Value = 0

$ws.Range("F4").Value = 734615892234.8064
$ws.Range("G4").Value = 88889835996.30263

$ws.Range("F5").Value = 558865056646.082
$ws.Range("G5").Value = 72600947639.16805

$ws.Range("F6").Value = 4579473077980.816
$ws.Range("G6").Value = 674619880691.7614

# Row 6 column E (M_POP) value was also corrected as part of this change.
$ws.Range("E6").Value = 2427884184.75
